# LogicComponentSequenceDiagram.pptx edit
# - Rename the ":Address" / "BookParser" actor label (two paragraphs) to a
#   single ":DiveLogParser" label (one paragraph, split across the existing
#   ":" + name runs).
# - Rename the "deletePerson(p)" call label to "deleteDiveSession(p)" and
#   let/force the textbox to grow to its auto-fit height.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {

    if ($shp.Id -eq 16) {
        # Was two paragraphs: ":Address" / "BookParser".
        # Becomes one paragraph: ":" + "DiveLogParser" runs.
        $tr = $shp.TextFrame.TextRange
        $tr.Text = ":DiveLogParser"

        # Re-touch the name portion in place so it stays its own run
        # (matching the original ":" + name run split).
        $tr2 = $shp.TextFrame.TextRange
        $nameLen = $tr2.Length - 1
        $c = $tr2.Characters(2, $nameLen)
        $c.Text = "DiveLogParser"
    }

    if ($shp.Id -eq 78) {
        # "deletePerson(p)" -> "deleteDiveSession(p)"; keep "(p)" run as-is.
        $tr = $shp.TextFrame.TextRange
        $c = $tr.Characters(1, 12)
        $c.Text = "deleteDiveSession"

        # This textbox auto-fits its height to the (now two-line) text;
        # pin it to the canonical auto-fit value.
        $shp.Height = 430887.0 / 12700.0
    }
}
